# "Beginning final month: March"
# Adds the 2021-03-05 (serial 44260) data row to Fallecido_Recuperado,
# opens the 2021-03-12 (serial 44267) placeholder date row on both sheets,
# and adds the full 2021-03-05 province breakdown to Provincias_Semanal.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Fallecido_Recuperado": fill in row 51 (date already present),
# then append row 52 with just the next reporting date.
# ---------------------------------------------------------------------
$wsFR = $wb.Worksheets.Item("Fallecido_Recuperado")

$wsFR.Cells.Item(51, 2).Value = 242660
$wsFR.Cells.Item(51, 3).Value = 3162
$wsFR.Cells.Item(51, 4).Value = 196484

# New row 52 - copy the date format from the row above, then set the value.
$wsFR.Range("A51").Copy()
$wsFR.Range("A52").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$wsFR.Cells.Item(52, 1).Value = 44267

# ---------------------------------------------------------------------
# Sheet "Provincias_Semanal": append the 32 province rows for 2021-03-05
# (serial 44260), then the lone next-date row (serial 44267).
# ---------------------------------------------------------------------
$wsPS = $wb.Worksheets.Item("Provincias_Semanal")

$provinces = @(
    @{ Name = "Distrito Nacional"; C = 6626.91; D = 496 },
    @{ Name = "Azua"; C = 1318.3; D = 32 },
    @{ Name = "Baoruco"; C = 1330.63; D = 13 },
    @{ Name = "Barahona"; C = 1377.02; D = 27 },
    @{ Name = "Dajabon"; C = 1501.5; D = 12 },
    @{ Name = "Duarte"; C = 2113.66; D = 200 },
    @{ Name = "Elias Pina"; C = 635.63; D = 6 },
    @{ Name = "El Seibo"; C = 931.47; D = 9 },
    @{ Name = "Espaillat"; C = 2118.87; D = 114 },
    @{ Name = "Independencia"; C = 1648.3; D = 11 },
    @{ Name = "La Altagracia"; C = 2544.9499999999998; D = 65 },
    @{ Name = "La Romana"; C = 2691.89; D = 134 },
    @{ Name = "La Vega"; C = 2402.9299999999998; D = 185 },
    @{ Name = "Maria Trinidad Sanchez"; C = 2153.84; D = 22 },
    @{ Name = "Monte Cristi"; C = 943.52; D = 23 },
    @{ Name = "Pedernales"; C = 1791.58; D = 4 },
    @{ Name = "Peravia"; C = 1019.07; D = 50 },
    @{ Name = "Puerto Plata"; C = 2255.2600000000002; D = 148 },
    @{ Name = "Hermanas Mirabal"; C = 2486.2199999999998; D = 37 },
    @{ Name = "Samana"; C = 778.35; D = 6 },
    @{ Name = "San Cristobal"; C = 1126.56; D = 145 },
    @{ Name = "San Juan"; C = 1621.75; D = 53 },
    @{ Name = "San Pedro de Macoris"; C = 1055.8599999999999; D = 49 },
    @{ Name = "Sanchez Ramirez"; C = 2025.95; D = 35 },
    @{ Name = "Santiago"; C = 2600.35; D = 496 },
    @{ Name = "Santiago Rodriguez"; C = 1955.79; D = 11 },
    @{ Name = "Valverde"; C = 1371.66; D = 41 },
    @{ Name = "Monsenor Nouel"; C = 1963.53; D = 47 },
    @{ Name = "Monte Plata"; C = 614.54999999999995; D = 31 },
    @{ Name = "Hato Mayor"; C = 907.16; D = 16 },
    @{ Name = "San Jose de Ocoa"; C = 1461.26; D = 15 },
    @{ Name = "Santo Domingo"; C = 1651.27; D = 629 }
)

$startRow = 1570
for ($i = 0; $i -lt $provinces.Count; $i++) {
    $row = $startRow + $i
    $prov = $provinces[$i]

    # Copy the date-formatted cell from the prior date block so the new
    # A-column cell reuses the existing date number format (style s="2").
    $wsPS.Range("A" + ($row - 32)).Copy()
    $wsPS.Range("A" + $row).PasteSpecial(-4122)
    $excel.CutCopyMode = $false

    $wsPS.Cells.Item($row, 1).Value = 44260
    $wsPS.Cells.Item($row, 2).Value = $prov.Name
    $wsPS.Cells.Item($row, 3).Value = $prov.C
    $wsPS.Cells.Item($row, 4).Value = $prov.D
}

$nextDateRow = $startRow + $provinces.Count
$wsPS.Range("A" + ($nextDateRow - 32)).Copy()
$wsPS.Range("A" + $nextDateRow).PasteSpecial(-4122)
$excel.CutCopyMode = $false
$wsPS.Cells.Item($nextDateRow, 1).Value = 44267

# ---------------------------------------------------------------------
# Restore selections/scroll so Provincias_Semanal remains the active tab
# (matches the workbook's stored activeTab), each sheet keeping its own
# last-used selection.
# ---------------------------------------------------------------------
$wsFR.Range("B52").Select() | Out-Null
$wsPS.Activate() | Out-Null
$wsPS.Range("B" + $nextDateRow).Select() | Out-Null
